$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 4.33
$ws.Range("K3").Value = 2.1
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.75
$ws.Range("AC3").Value = 9
$ws.Range("AG3").Value = 451
$ws.Range("AO3").Value = 9.5
$ws.Range("AT3").Value = 2.63
$ws.Range("AU3").Value = 8.5
$ws.Range("AW3").Value = 6
$ws.Range("G5").Value = 1.53
$ws.Range("I5").Value = 7
$ws.Range("K5").Value = 2.05
$ws.Range("M5").Value = 1.11
$ws.Range("N5").Value = 6.5
$ws.Range("X5").Value = 5.5
$ws.Range("O7").Value = 1.5
$ws.Range("P7").Value = 2.5
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("H8").Value = 3.75
$ws.Range("J8").Value = 2.87
$ws.Range("K8").Value = 2.35
$ws.Range("L8").Value = 2.87
$ws.Range("Q8").Value = 1.5
$ws.Range("R8").Value = 2.27
$ws.Range("U8").Value = 1.45
$ws.Range("V8").Value = 2.37
$ws.Range("W8").Value = 12.5
$ws.Range("X8").Value = 15.5
$ws.Range("Z8").Value = 28
$ws.Range("AC8").Value = 16.5
$ws.Range("AD8").Value = 7.8
$ws.Range("AE8").Value = 11.5
$ws.Range("AF8").Value = 37
$ws.Range("AH8").Value = 12
$ws.Range("AI8").Value = 15
$ws.Range("AJ8").Value = 9.75
$ws.Range("AL8").Value = 17
$ws.Range("AM8").Value = 21
$ws.Range("AN8").Value = 4.75
$ws.Range("AP8").Value = 16
$ws.Range("AT8").Value = 3.1
$ws.Range("AU8").Value = 6.3
$ws.Range("AW8").Value = 4.75
$ws.Range("AY8").Value = 16
$ws.Range("BA8").Value = 60
$ws.Range("BC8").Value = 500
$ws.Range("G14").Value = 2.15
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 2.77
$ws.Range("K14").Value = 2.12
$ws.Range("L14").Value = 3.6
$ws.Range("N14").Value = 6.9
$ws.Range("P14").Value = 2.95
$ws.Range("U14").Value = 1.85
$ws.Range("V14").Value = 1.85
$ws.Range("W14").Value = 7.1
$ws.Range("X14").Value = 10
$ws.Range("Y14").Value = 9
$ws.Range("Z14").Value = 20
$ws.Range("AA14").Value = 18.5
$ws.Range("AC14").Value = 6.9
$ws.Range("AD14").Value = 6.5
$ws.Range("AE14").Value = 15.5
$ws.Range("AF14").Value = 80
$ws.Range("AH14").Value = 8.75
$ws.Range("AI14").Value = 15
$ws.Range("AJ14").Value = 11
$ws.Range("AK14").Value = 37
$ws.Range("AL14").Value = 28
$ws.Range("AM14").Value = 40
$ws.Range("AN14").Value = 4.05
$ws.Range("AO14").Value = 11.25
$ws.Range("AQ14").Value = 45
$ws.Range("AU14").Value = 7.5
$ws.Range("AW14").Value = 4.9
$ws.Range("AX14").Value = 16.5
$ws.Range("AY14").Value = 26
$ws.Range("AZ14").Value = 80
